# Insert 4 new "data quality assertion" rows (duplicate_value_combination:*
# and geo_spatial_accuracy_precision:*) ahead of the existing
# "date_format_validation:*" block, pushing all subsequent rows down by 4
# (old row 28 -> new row 32, ... old row 48 -> new row 52) in both the
# "Data quality assertion" sheet and the "Definition of assertions" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Data quality assertion" (2 columns: A = assertion key, B = blank
# "New Use Case Name" placeholder column)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data quality assertion")

$ws2.Range("A28:B31").EntireRow.Insert()

$ws2.Cells.Item(28, 1).Value = "duplicate_value_combination:inferred_duplicate"
$ws2.Cells.Item(28, 2).Value = ""
$ws2.Cells.Item(29, 1).Value = "duplicate_value_combination:inferred_nonduplicate"
$ws2.Cells.Item(29, 2).Value = ""
$ws2.Cells.Item(30, 1).Value = "geo_spatial_accuracy_precision:low_precision"
$ws2.Cells.Item(30, 2).Value = ""
$ws2.Cells.Item(31, 1).Value = "geo_spatial_accuracy_precision:high_precision"
$ws2.Cells.Item(31, 2).Value = ""

# ---------------------------------------------------------------------
# Sheet "Definition of assertions" (6 columns: A = assertion key,
# B = Category, C = Input field (RDF), D = Label, E = Simple Rule
# Definition, F = Expanded Rule Definition)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Definition of assertions")

$ws3.Range("A28:F31").EntireRow.Insert()

$ws3.Cells.Item(28, 1).Value = "duplicate_value_combination:inferred_duplicate"
$ws3.Cells.Item(28, 2).Value = "data_quality"
$ws3.Cells.Item(28, 3).Value = "various_fields"
$ws3.Cells.Item(28, 4).Value = "inferred_duplicate"
$ws3.Cells.Item(28, 5).Value = "Indicates that the record has a combination of values across multiple fields that are identical to other records."
$ws3.Cells.Item(28, 6).Value = "If the record has a combination of values across the specified fields that is identical to other records in the dataset, label it as 'duplicate_combination'. This implies redundancy in data values for multiple records."

$ws3.Cells.Item(29, 1).Value = "duplicate_value_combination:inferred_nonduplicate"
$ws3.Cells.Item(29, 2).Value = "data_quality"
$ws3.Cells.Item(29, 3).Value = "various_fields"
$ws3.Cells.Item(29, 4).Value = "inferred_nonduplicate"
$ws3.Cells.Item(29, 5).Value = "Indicates that the record has a unique combination of values across multiple fields that is not shared by other records."
$ws3.Cells.Item(29, 6).Value = "If the record has a unique combination of values across the specified fields, label it as 'unique_combination'. This means that no other records share this exact combination."

$ws3.Cells.Item(30, 1).Value = "geo_spatial_accuracy_precision:low_precision"
$ws3.Cells.Item(30, 2).Value = "geo"
$ws3.Cells.Item(30, 3).Value = "geo:hasMetricSpatialAccuracy"
$ws3.Cells.Item(30, 4).Value = "low_precision"
$ws3.Cells.Item(30, 5).Value = "Indicates that the spatial accuracy is low, either because the value of coordinateUncertaintyInMeters is empty or exceeds 10,000 meters."
$ws3.Cells.Item(30, 6).Value = "If the 'coordinateUncertaintyInMeters' field is empty or its value exceeds 10,000, label the record as 'low_precision'. This indicates that the precision of the spatial accuracy is insufficient."

$ws3.Cells.Item(31, 1).Value = "geo_spatial_accuracy_precision:high_precision"
$ws3.Cells.Item(31, 2).Value = "geo"
$ws3.Cells.Item(31, 3).Value = "geo:hasMetricSpatialAccuracy"
$ws3.Cells.Item(31, 4).Value = "high_precision"
$ws3.Cells.Item(31, 5).Value = "Indicates that the spatial accuracy is high, meaning the value of coordinateUncertaintyInMeters is less than or equal to 10,000 meters."
$ws3.Cells.Item(31, 6).Value = "If the 'coordinateUncertaintyInMeters' field contains a value of 10,000 or less, label the record as 'high_precision'. This indicates that the precision of the spatial accuracy is adequate."
